# Update "all_Servos_angles" workbook: add a new "waist" body-part section
# to the "2019" sheet, and highlight the related head/neck rotation rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuille 1")
$ws2 = $wb.Worksheets.Item("2019")

# Green fill colour used elsewhere in the sheet for "section" rows
# (FF92D050 in RGB -> BGR packed value for the COM Color property).
$green = 5296274

# --- "2019" sheet: insert two new rows after the "Head - Neck" block ----
# Row 8 becomes a blank spacer (matching the style used throughout the
# sheet) and row 9 introduces the new "waist" servo entry.
$ws2.Rows("8:9").Insert()

$ws2.Cells.Item(9, 1).Value = "waist"
$ws2.Cells.Item(9, 3).Value = 45
$ws2.Cells.Item(9, 4).Value = 135
$ws2.Cells.Item(9, 6).Value = 90

# Highlight the "waist" label plus the related "rothead" / "neck" rows so
# they stand out as the newly documented body parts.
$ws2.Cells.Item(9, 1).Interior.Color = $green
$ws2.Cells.Item(6, 1).Interior.Color = $green
$ws2.Cells.Item(7, 1).Interior.Color = $green

# --- Refresh the saved cursor/scroll position on both sheets -----------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("E14").Select()

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("C9").Select()
